# Update LeetCode Hot 100
#
# Nudges four of the small Consolas label textboxes in the "dp[i] formula"
# diagram horizontally (their vertical position / size is untouched):
#   - "i"              (shape id 46): x 6252337 -> 6223762 EMU
#   - "dp[i-1]"        (shape id 54): x 5609212 -> 5533012 EMU
#   - "i-dp[i-1]-2"     (shape id 55): x 3673182 -> 3720807 EMU
#   - "i-dp[i-1]-1"     (shape id 64): x 4539022 -> 4491397 EMU
#
# Shape.Left/Top on the PowerPoint object model are expressed in points
# (1 pt = 12700 EMU), so the EMU deltas from the source diagram are
# converted to points below. The literals carry enough decimal digits to
# land on the exact target EMU after the host's internal point<->EMU
# conversion.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

foreach ($sh in $s.Shapes) {
    switch ($sh.Id) {
        46 { $sh.Left = 490.0600280762 }   # "i"            -> x=6223762
        54 { $sh.Left = 435.6702575684 }   # "dp[i-1]"      -> x=5533012
        55 { $sh.Left = 292.9769592285 }   # "i-dp[i-1]-2"  -> x=3720807
        64 { $sh.Left = 353.6533203125 }   # "i-dp[i-1]-1"  -> x=4491397
    }
}
